$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Treinamento")
$ws.Activate()

$ws.Range("B7").Value = 2
$ws.Range("B21").Value = 2
$ws.Range("B32").Value = 2
$ws.Range("B34").Value = 2
$ws.Range("B39").Value = 2
$ws.Range("B40").Value = 2
$ws.Range("B57").Value = 2
$ws.Range("B58").Value = 2
$ws.Range("B63").Value = 2
$ws.Range("B69").Value = 2
$ws.Range("B87").Value = 2
$ws.Range("B114").Value = 2
$ws.Range("B122").Value = 2
$ws.Range("B127").Value = 2
$ws.Range("B161").Value = 2
$ws.Range("B179").Value = 2
$ws.Range("B190").Value = 2
$ws.Range("B192").Value = 2
$ws.Range("B194").Value = 2
$ws.Range("B209").Value = 2
$ws.Range("B227").Value = 2
$ws.Range("B239").Value = 2
$ws.Range("B243").Value = 2
$ws.Range("B252").Value = 2
$ws.Range("B265").Value = 3
$ws.Range("B268").Value = 2
$ws.Range("B275").Value = 2
$ws.Range("B295").Value = 2
$ws.Range("B297").Value = 2

$ws.Range("T297").Select()
